# Apply crypto price/volume updates per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is plain (non-numeric-looking) text -- direct assignment keeps them as text
$ws.Range("D2").Value = "26.714.01"
$ws.Range("E2").Value = "  -2.46%  "
$ws.Range("D3").Value = "1.556.81"
$ws.Range("E3").Value = "  -0.50%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("E6").Value = "  -1.82%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.60%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "1.776.87"
$ws.Range("E12").Value = "  -0.55%  "
$ws.Range("D13").Value = "1.558.31"
$ws.Range("E13").Value = "  -0.57%  "
$ws.Range("E14").Value = "  -2.34%  "
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("E16").Value = "  -2.77%  "
$ws.Range("D17").Value = "26.737.08"
$ws.Range("E17").Value = "  -2.42%  "
$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("E18").Value = "  +1.27%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("E19").Value = "  +0.63%  "
$ws.Range("D20").Value = "0.0₃0675"
$ws.Range("E20").Value = "  -1.93%  "
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  -0.48%  "
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("E25").Value = "  -0.50%  "
$ws.Range("E26").Value = "  +0.87%  "
$ws.Range("E27").Value = "  -1.07%  "
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  -1.14%  "
$ws.Range("E30").Value = "  -1.37%  "
$ws.Range("E31").Value = "  -3.99%  "
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").Value = "1.383.59"
$ws.Range("E33").Value = "  +1.51%  "
$ws.Range("E34").Value = "  -1.27%  "
$ws.Range("E35").Value = "  +0.96%  "
$ws.Range("E36").Value = "  -0.65%  "
$ws.Range("E37").Value = "  -4.19%  "
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("E39").Value = "  -2.46%  "
$ws.Range("E40").Value = "  -1.27%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +1.93%  "
$ws.Range("E43").Value = "  +2.05%  "
$ws.Range("E44").Value = "  +1.48%  "
$ws.Range("E45").Value = "  -1.31%  "
$ws.Range("E46").Value = "  -1.46%  "
$ws.Range("D47").Value = "1.690.36"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").Value = "0.0₇0984"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("E50").Value = "  -0.27%  "
$ws.Range("E51").Value = "  -1.06%  "

# Cells whose new value looks like a number (e.g. "205.85") but must stay stored as text,
# matching the original inline-string cell type. Force text format before assigning so Excel
# does not auto-convert the literal into a numeric value.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "205.85"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.490"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "21.85"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0582"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0862"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.73"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.511"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "61.54"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.35"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.14"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.08"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.01"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.56"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "6.77"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.82"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.15"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.91"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.28"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.933"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.518"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.811"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.992"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.39"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.18"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.05"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.49"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0492"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0945"

# Remove the temporary text-number-format override so the cells end up with no explicit
# style again (matching the original workbook, which left these cells unstyled).
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D10").ClearFormats()
$ws.Range("D11").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D24").ClearFormats()
$ws.Range("D25").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D32").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D36").ClearFormats()
$ws.Range("D37").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D42").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D44").ClearFormats()
$ws.Range("D46").ClearFormats()
$ws.Range("D48").ClearFormats()
$ws.Range("D50").ClearFormats()
$ws.Range("D51").ClearFormats()
